$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("E2").Value = "2016-03-22 07:09:29"
$wsZh.Range("E3").Value = "2016-03-22 07:09:29"
$wsZh.Range("H2").Value = "2016-03-22 07:10:16"
$wsZh.Range("H3").Value = "2016-03-22 07:10:16"

$wsDe.Range("E2").Value = "2016-03-22 07:09:38"
$wsDe.Range("E3").Value = "2016-03-22 07:09:38"
$wsDe.Range("H2").Value = "2016-03-22 07:10:35"
$wsDe.Range("H3").Value = "2016-03-22 07:10:35"
